$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and symbol-list (column E) updates scraped on 2022-12-12.
# Values are kept as TEXT (matching the original inlineStr/text storage of this
# sheet) by forcing a "@" (Text) number format before assignment, then resetting
# the cell style back to "Normal" so no stray formatting is left behind.
$updates = @{
    "D2" = "276.27"
    "D3" = "20.83"
    "D4" = "6.211"
    "D5" = "0.06168"
    "D6" = "3.579"
    "D7" = "6.575"
    "D8" = "1.504"
    "D9" = "0.8183"
    "D10" = "0.01382"
    "D11" = "0.1620"
    "D12" = "0.08313"
    "D13" = "0.03681"
    "D14" = "0.03145"
    "D15" = "0.09118"
    "D16" = "3.708"
    "D17" = "0.001633"
    "D18" = "0.04673"
    "D19" = "0.006427"
    "D20" = "0.006162"
    "D21" = "0.001069"
    "D22" = "0.0001501"
    "D23" = "3.775"
    "D24" = "2.231"
    "D25" = "0.3388"
    "D26" = "0.1222"
    "E27" = "26AAXTokenAAB"
    "D40" = "0.04661"
    "D41" = "0.007029"
    "D42" = "0.1100"
    "D43" = "0.003522"
    "D44" = "0.01114"
    "D45" = "0.00006450"
    "D46" = "0.00000000751"
    "D47" = "0.8416"
    "D48" = "0.002720"
    "D49" = "0.00001902"
    "E49" = "48CryptobidCoinCBCBestin24h"
    "D50" = "0.01242"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}

Write-Output ("Updated " + $updates.Count + " cells")
